$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy header style (bold font, thin border, center/top alignment) from X1
#     onto the newly added header cells Y1:AF1 (format-only paste; value set after) ---
$ws.Range("X1").Copy() | Out-Null
$ws.Range("Y1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- New resequencing coverage-abundance data (columns B..AF for rows 1..7) ---
# Row 1 holds the numeric column headers (0..30); rows 2-7 hold the per-sample
# abundance fractions. Trailing $null entries are left blank (no cell written),
# matching the source data which has fewer samples for those reaction pairs.
$rowData = @{
    1 = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30)
    2 = @(0.632435353200033, 0.4393974168277774, 0.4493089066013006, 0.6335194188791039, 0.6673157768427789, 0.6137177741528178, 0.5625848866572797, 0.581521570053209, 0.5464768148922564, 0.3293768863682892, 0.4096584823401687, 0.2864167677962447, 0.3382050257764608, 0.5418061185836966, 0.3022745271344543, 0.2649012839886506, 0.5853761838818446, 0.5488430109465758, 0.3845235247633266, 0.3484211687312401, 0.4487883201861856, 0.5372157293234545, 0.5263266359644805, 0.5461027859213058, 0.4871100146020673, 0.5653076881099184, 0.5038166212121797, 0.5844982878134367, 0.6198687257631004, 0.55219714399938, 0.5606273828260427)
    3 = @(0.3675646467999671, 0.5606025831722226, 0.5506910933986994, 0.3664805811208961, 0.3326842231572211, 0.3862822258471822, 0.4374151133427203, 0.418478429946791, 0.4535231851077436, 0.6706231136317108, 0.5903415176598313, 0.7135832322037553, 0.6617949742235391, 0.4581938814163033, 0.6977254728655456, 0.7350987160113494, 0.4146238161181554, 0.4511569890534242, 0.6154764752366734, 0.6515788312687599, 0.5512116798138145, 0.4627842706765454, 0.4736733640355195, 0.4538972140786942, 0.5128899853979327, 0.4346923118900817, 0.4961833787878202, 0.4155017121865633, 0.3801312742368996, 0.4478028560006199, 0.4393726171739572)
    4 = @(0.5933141162398911, 0.6903071402079209, 0.5549258000377028, 0.6379465615468691, 0.6385690069986263, 0.5443587682023839, 0.6993821024193529, 0.5546314308031363, 0.714871391653405, 0.5198485053999897, 0.6281848790015272, 0.5479949405358678, 0.5402916407514635, 0.5847128666817604, 0.5421145898894826, 0.6731729945773508, 0.5314466969516085, 0.7208905243865437, 0.7592157288413576, 0.5962721292230511, 0.661048092863724, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    5 = @(0.4066858837601089, 0.3096928597920791, 0.4450741999622973, 0.3620534384531308, 0.3614309930013738, 0.4556412317976161, 0.3006178975806471, 0.4453685691968637, 0.285128608346595, 0.4801514946000104, 0.3718151209984728, 0.4520050594641323, 0.4597083592485365, 0.4152871333182396, 0.4578854101105175, 0.3268270054226491, 0.4685533030483915, 0.2791094756134563, 0.2407842711586423, 0.403727870776949, 0.338951907136276, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    6 = @(0.3086710021568247, 0.246898793886378, 0.147097704055629, 0.1629162617895906, 0.1930483433629932, 0.2081105927477706, 0.338793300765507, 0.2953945009440818, 0.2565166105438769, 0.2611125379647902, 0.142011226138745, 0.1337598637833575, 0.1031692489390861, 0.3345501341576917, 0.1847636367775405, 0.1215432442380801, 0.0787665057798845, 0.1353623430768313, 0.1284474604153815, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    7 = @(0.6913289978431753, 0.753101206113622, 0.8529022959443711, 0.8370837382104094, 0.8069516566370067, 0.7918894072522294, 0.661206699234493, 0.7046054990559182, 0.7434833894561232, 0.7388874620352098, 0.857988773861255, 0.8662401362166425, 0.8968307510609139, 0.6654498658423083, 0.8152363632224595, 0.8784567557619198, 0.9212334942201156, 0.8646376569231687, 0.8715525395846185, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item([int]$r, $i + 2).Value = $v
        }
    }
}

# --- Cells that were already blank (empty, text-typed placeholder) in the source
#     workbook and stay blank here too. Explicitly re-clear them: a plain load/save
#     round-trip otherwise leaks shared-string #0 into these empty cells. ---
$reblankCells = @("W4", "X4", "W5", "X5", "U6", "V6", "W6", "X6", "U7", "V7", "W7", "X7")
foreach ($addr in $reblankCells) {
    $ws.Range($addr).ClearContents()
}

Write-Output "done"
